$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.756.96"
Set-TextValue $ws.Range("E2") "  +5.42%  "

Set-TextValue $ws.Range("D3") "2.277.80"
Set-TextValue $ws.Range("E3") "  +3.71%  "

Set-TextValue $ws.Range("E4") "  +0.18%  "

Set-TextValue $ws.Range("D5") "232.84"
Set-TextValue $ws.Range("E5") "  +1.64%  "

Set-TextValue $ws.Range("D6") "0.638"
Set-TextValue $ws.Range("E6") "  +3.58%  "

Set-TextValue $ws.Range("D7") "65.73"
Set-TextValue $ws.Range("E7") "  +9.41%  "

Set-TextValue $ws.Range("E8") "  +0.12%  "

Set-TextValue $ws.Range("D9") "0.430"
Set-TextValue $ws.Range("E9") "  +7.39%  "

Set-TextValue $ws.Range("E10") "  +17.42%  "

Set-TextValue $ws.Range("D11") "57.60"
Set-TextValue $ws.Range("E11") "  +1.30%  "

Set-TextValue $ws.Range("D12") "26.35"
Set-TextValue $ws.Range("E12") "  +18.88%  "

Set-TextValue $ws.Range("E13") "  +0.00%  "

Set-TextValue $ws.Range("D14") "2.616.97"
Set-TextValue $ws.Range("E14") "  +3.79%  "

Set-TextValue $ws.Range("D15") "16.00"
Set-TextValue $ws.Range("E15") "  +4.36%  "

Set-TextValue $ws.Range("E16") "  +6.03%  "

Set-TextValue $ws.Range("D17") "0.832"
Set-TextValue $ws.Range("E17") "  +5.34%  "

Set-TextValue $ws.Range("D18") "2.275.47"
Set-TextValue $ws.Range("E18") "  +3.74%  "

Set-TextValue $ws.Range("D19") "43.767.39"
Set-TextValue $ws.Range("E19") "  +5.68%  "

Set-TextValue $ws.Range("E20") "  +12.37%  "

Set-TextValue $ws.Range("D21") "74.18"
Set-TextValue $ws.Range("E21") "  +3.25%  "

Set-TextValue $ws.Range("D22") "6.13"
Set-TextValue $ws.Range("E22") "  +1.83%  "

Set-TextValue $ws.Range("D23") "251.14"
Set-TextValue $ws.Range("E23") "  +3.85%  "

Set-TextValue $ws.Range("E24") "  +0.12%  "

Set-TextValue $ws.Range("D25") "2.49"
Set-TextValue $ws.Range("E25") "  +6.23%  "

Set-TextValue $ws.Range("D26") "2.33"
Set-TextValue $ws.Range("E26") "  +1.80%  "

Set-TextValue $ws.Range("D27") "10.14"
Set-TextValue $ws.Range("E27") "  +5.98%  "

Set-TextValue $ws.Range("D28") "173.44"
Set-TextValue $ws.Range("E28") "  +2.83%  "

Set-TextValue $ws.Range("D29") "21.01"
Set-TextValue $ws.Range("E29") "  +6.81%  "

Set-TextValue $ws.Range("E30") "  -1.21%  "

Set-TextValue $ws.Range("E31") "  +0.27%  "

Set-TextValue $ws.Range("D32") "2.78"
Set-TextValue $ws.Range("E32") "  +8.53%  "

Set-TextValue $ws.Range("D33") "0.124"
Set-TextValue $ws.Range("E33") "  +3.45%  "

Set-TextValue $ws.Range("E34") "  +6.96%  "

Set-TextValue $ws.Range("D35") "5.07"
Set-TextValue $ws.Range("E35") "  +2.17%  "

Set-TextValue $ws.Range("D36") "4.77"
Set-TextValue $ws.Range("E36") "  +3.74%  "

Set-TextValue $ws.Range("D37") "3.84"
Set-TextValue $ws.Range("E37") "  +8.96%  "

Set-TextValue $ws.Range("D38") "6.76"
Set-TextValue $ws.Range("E38") "  +7.48%  "

Set-TextValue $ws.Range("D39") "2.36"
Set-TextValue $ws.Range("E39") "  +0.77%  "

Set-TextValue $ws.Range("E40") "  +6.00%  "

Set-TextValue $ws.Range("E41") "  +0.26%  "

Set-TextValue $ws.Range("D42") "17.76"
Set-TextValue $ws.Range("E42") "  +9.00%  "

Set-TextValue $ws.Range("B43") "FTXToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D43") "4.60"
Set-TextValue $ws.Range("E43") "  +6.47%  "

Set-TextValue $ws.Range("B44") "FraxShare"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "8.41"
Set-TextValue $ws.Range("E44") "  -1.02%  "

Set-TextValue $ws.Range("D45") "10.56"
Set-TextValue $ws.Range("E45") "  +23.85%  "

Set-TextValue $ws.Range("D46") "0.0975"
Set-TextValue $ws.Range("E46") "  +2.89%  "

Set-TextValue $ws.Range("B47") "Aave"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "98.40"
Set-TextValue $ws.Range("E47") "  +1.82%  "

Set-TextValue $ws.Range("B48") "TrustWalletToken"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D48") "1.21"
Set-TextValue $ws.Range("E48") "  +1.55%  "

Set-TextValue $ws.Range("D49") "1.481.36"
Set-TextValue $ws.Range("E49") "  +1.48%  "

Set-TextValue $ws.Range("E50") "  +6.48%  "

Set-TextValue $ws.Range("D51") "0.000205"
Set-TextValue $ws.Range("E51") "  -14.23%  "
